# Fix the CODEVALUE shared string "Koodisto 6000" -> "Koodisto6000"
# (removing the invalid embedded space so it is a valid codeValue).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CodeSchemes")
$ws.Range("A2").Value = "Koodisto6000"
